$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '76.402.69'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '

$ws.Cells.Item(3, 4).Value = '3.077.05'
$ws.Cells.Item(3, 5).Value = '  +4.62%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '''198.36'
$ws.Cells.Item(5, 5).Value = '  +0.03%  '

$ws.Cells.Item(6, 4).Value = '''619.07'
$ws.Cells.Item(6, 5).Value = '  +4.14%  '

$ws.Cells.Item(7, 5).Value = '  +0.02%  '

$ws.Cells.Item(8, 2).Value = 'Dogecoin'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(8, 4).Value = '''0.215'
$ws.Cells.Item(8, 5).Value = '  +7.83%  '

$ws.Cells.Item(9, 2).Value = 'XRP'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(9, 4).Value = '''0.552'
$ws.Cells.Item(9, 5).Value = '  +0.51%  '

$ws.Cells.Item(10, 4).Value = '3.074.53'
$ws.Cells.Item(10, 5).Value = '  +4.65%  '

$ws.Cells.Item(11, 5).Value = '  +0.60%  '

$ws.Cells.Item(12, 5).Value = '  -0.02%  '

$ws.Cells.Item(13, 4).Value = '''5.23'
$ws.Cells.Item(13, 5).Value = '  +7.05%  '

$ws.Cells.Item(14, 4).Value = '3.644.47'
$ws.Cells.Item(14, 5).Value = '  +4.76%  '

$ws.Cells.Item(15, 4).Value = '''29.29'
$ws.Cells.Item(15, 5).Value = '  +2.91%  '

$ws.Cells.Item(16, 4).Value = '''0.0000197'
$ws.Cells.Item(16, 5).Value = '  +3.85%  '

$ws.Cells.Item(17, 4).Value = '76.265.40'
$ws.Cells.Item(17, 5).Value = '  -0.24%  '

$ws.Cells.Item(18, 4).Value = '3.068.52'
$ws.Cells.Item(18, 5).Value = '  +4.74%  '

$ws.Cells.Item(19, 5).Value = '  -0.37%  '

$ws.Cells.Item(20, 4).Value = '''9.01'
$ws.Cells.Item(20, 5).Value = '  +3.04%  '

$ws.Cells.Item(21, 4).Value = '''384.66'
$ws.Cells.Item(21, 5).Value = '  +2.68%  '

$ws.Cells.Item(22, 5).Value = '  +13.86%  '

$ws.Cells.Item(23, 4).Value = '''4.51'
$ws.Cells.Item(23, 5).Value = '  +4.20%  '

$ws.Cells.Item(24, 5).Value = '  +0.79%  '

$ws.Cells.Item(25, 4).Value = '''4.58'
$ws.Cells.Item(25, 5).Value = '  +7.46%  '

$ws.Cells.Item(26, 4).Value = '3.231.39'
$ws.Cells.Item(26, 5).Value = '  +4.70%  '

$ws.Cells.Item(27, 4).Value = '''72.44'

$ws.Cells.Item(28, 5).Value = '  +0.11%  '

$ws.Cells.Item(29, 5).Value = '  +4.24%  '

$ws.Cells.Item(30, 5).Value = '  +0.46%  '

$ws.Cells.Item(31, 4).Value = '''0.994'
$ws.Cells.Item(31, 5).Value = '  -0.59%  '

$ws.Cells.Item(32, 5).Value = '  -0.25%  '

$ws.Cells.Item(33, 5).Value = '  +3.84%  '

$ws.Cells.Item(34, 4).Value = '''501.30'
$ws.Cells.Item(34, 5).Value = '  +0.51%  '

$ws.Cells.Item(35, 5).Value = '  +6.01%  '

$ws.Cells.Item(36, 4).Value = '''0.130'
$ws.Cells.Item(36, 5).Value = '  +17.11%  '

$ws.Cells.Item(37, 5).Value = '  +0.02%  '

$ws.Cells.Item(38, 4).Value = '''20.89'
$ws.Cells.Item(38, 5).Value = '  +3.71%  '

$ws.Cells.Item(39, 4).Value = '''163.40'
$ws.Cells.Item(39, 5).Value = '  -1.19%  '

$ws.Cells.Item(40, 4).Value = '''194.96'
$ws.Cells.Item(40, 5).Value = '  +8.55%  '

$ws.Cells.Item(41, 5).Value = '  +0.54%  '

$ws.Cells.Item(42, 5).Value = '  -3.53%  '

$ws.Cells.Item(43, 4).Value = '''0.102'
$ws.Cells.Item(43, 5).Value = '  -7.04%  '

$ws.Cells.Item(44, 5).Value = '  +0.05%  '

$ws.Cells.Item(45, 4).Value = '''0.798'
$ws.Cells.Item(45, 5).Value = '  +19.95%  '

$ws.Cells.Item(46, 4).Value = '''5.17'
$ws.Cells.Item(46, 5).Value = '  +5.34%  '

$ws.Cells.Item(47, 5).Value = '  +6.85%  '

$ws.Cells.Item(48, 5).Value = '  +1.22%  '

$ws.Cells.Item(49, 2).Value = 'OKB'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(49, 4).Value = '''40.98'
$ws.Cells.Item(49, 5).Value = '  +2.55%  '

$ws.Cells.Item(50, 2).Value = 'dogwifhat'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(50, 4).Value = '''2.45'
$ws.Cells.Item(50, 5).Value = '  +5.18%  '

$ws.Cells.Item(51, 4).Value = '''0.599'
$ws.Cells.Item(51, 5).Value = '  +0.61%  '

